# Applies the "Added many more features" edit to the Montezuma's Treasure
# slot game review document: updates the title, the "what we like" /
# "what we don't like" bullet lists, and the closing bold title + italic
# meta description paragraph.

$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title (appears twice: the H1 heading and the bold line near the end)
Replace-All "Play Montezuma's Treasure for Free - Slot Game Review" `
            "Play Montezuma’s Treasure for Free - Exciting Aztec Slot Game"

# "What we like" bullets
Replace-All "Visually impressive design and theme" `
            "Visually impressive graphics, effects, and sounds"

Replace-All "Engaging gameplay for experienced gamblers" `
            "Engaging gameplay with double spins and ten paylines"

Replace-All "High-quality graphics and sound effects" `
            "Exciting bonus game with multiple levels and prizes"

Replace-All "Exciting bonuses and opportunities" `
            "High-quality graphics, sounds, and effects for an immersive experience"

# "What we don't like" bullets
Replace-All "Limited paylines available" `
            "Limited selection of betting options with minimum bet at 0.25 cents"

Replace-All "Minimum bet may be too high for some players" `
            "Not suitable for players not interested in historical themes"

# Closing italic meta-description paragraph
Replace-All "Explore the ancient Aztec civilization with Montezuma's Treasure. Enjoy engaging gameplay, exciting bonuses, and high-quality graphics for free." `
            "Read our review of Montezuma’s Treasure and play this exciting Aztec-themed slot game for free. Experience the adventure today!"
